# Updated symbol list on Tue Feb 14 03:43:42 UTC 2023 with GitHub Actions
# Applies refreshed Price (D) and Volume(1h) (E) values to the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Address, $Text) {
    $cell = $Sheet.Range($Address)
    # Force the cell to be treated as literal text so values like "291.67"
    # or "-6.16%" are not reinterpreted by Excel as numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    # Restore the default (unstyled) cell style, matching the source file
    # which carries no explicit style on these data cells.
    $cell.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  D = "291.67";     E = "-6.16%" }
    @{ Row = 3;  D = "40.49";      E = "1.47%" }
    @{ Row = 4;  D = "5.019";      E = "-1.63%" }
    @{ Row = 5;  D = "0.07340";    E = "-3.09%" }
    @{ Row = 6;  D = "4.295";      E = "-0.13%" }
    @{ Row = 7;  E = "-7.20%" }
    @{ Row = 8;  D = "0.9222";     E = "-0.78%" }
    @{ Row = 10; D = "0.1219";     E = "0.44%" }
    @{ Row = 11; D = "0.1725";     E = "-4.91%" }
    @{ Row = 12; D = "0.08606";    E = "-4.61%" }
    @{ Row = 13; D = "0.04270";    E = "2.73%" }
    @{ Row = 14; D = "0.1053";     E = "-0.09%" }
    @{ Row = 15; D = "0.001280";   E = "-0.26%" }
    @{ Row = 16; D = "0.005780";   E = "-2.75%" }
    @{ Row = 17; E = "-0.36%" }
    @{ Row = 18; D = "0.3287";     E = "-2.01%" }
    @{ Row = 19; D = "7.704";      E = "0.61%" }
    @{ Row = 20; D = "0.1390";     E = "2.83%" }
    @{ Row = 21; D = "0.2748";     E = "-2.25%" }
    @{ Row = 22; D = "0.03932";    E = "-2.27%" }
    @{ Row = 23; D = "0.001261";   E = "-0.47%" }
    @{ Row = 24; D = "0.003777";   E = "-7.36%" }
    @{ Row = 25; D = "0.0001282";  E = "0.84%" }
    @{ Row = 26; D = "0.0003726";  E = "-95.05%" }
    @{ Row = 38; D = "0.02301";    E = "-5.04%" }
    @{ Row = 39; D = "0.04973";    E = "-3.28%" }
    @{ Row = 41; D = "0.007702";   E = "-0.41%" }
    @{ Row = 42; E = "-1.20%" }
    @{ Row = 43; D = "0.007363";   E = "-3.87%" }
    @{ Row = 44; D = "0.007788";   E = "-3.49%" }
    @{ Row = 45; D = "0.3171";     E = "2.06%" }
    @{ Row = 46; D = "0.00006346"; E = "-3.88%" }
    @{ Row = 47; E = "0.06%" }
    @{ Row = 48; D = "0.02116";    E = "-91.86%" }
    @{ Row = 49; E = "0.06%" }
    @{ Row = 50; E = "0.06%" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $addr = "D" + $u.Row
        Set-TextValue $ws $addr $u.D
    }
    if ($u.ContainsKey("E")) {
        $addr = "E" + $u.Row
        Set-TextValue $ws $addr $u.E
    }
}
